# Update the caption text box on slide 1 ("ROC and Precision-Recall Curves")
# The second paragraph is rewritten:
#   Old: "Their computation from the contingency table and shape when evaluating
#         data-driven models trained with balanced and imbalanced training datasets"
#   New: "Computation from contingency tables and performance characteristics under
#         balanced and imbalanced datasets"
# The new text keeps two runs: the first (non-underlined) run carries the bulk of the
# sentence, the second (also non-underlined) run carries "imbalanced datasets" -
# the underlined "computation"/"shape" emphasis from the old text is dropped.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "TextBox 861") {
        $targetShape = $candidate
        break
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the run boundaries (1-indexed) of the second paragraph by searching the
# full text, so the script is resilient to minor offset differences.
$full = $tr.Text
$p2Start = $full.IndexOf("Their ") + 1
$oldRun4Start = $full.IndexOf("shape when evaluating") + 1
$oldRun4Len = 5
$oldRun3Start = $full.IndexOf(" from the contingency table and ") + 1
$oldRun3Len = 33
$oldRun2Start = $full.IndexOf("computation") + 1
$oldRun2Len = 11

# Delete (right-to-left, so earlier offsets stay valid) the middle runs that carried
# the underlined emphasis ("computation" ... "shape") together with the connecting
# text between them - everything between "Their " and the final trailing clause.
$tr.Characters($oldRun4Start, $oldRun4Len).Text = ""
$tr.Characters($oldRun3Start, $oldRun3Len).Text = ""
$tr.Characters($oldRun2Start, $oldRun2Len).Text = ""

# What remains in paragraph 2 is now: "Their " + " when evaluating data-driven
# models trained with balanced and imbalanced training datasets" (both non-underlined).
$run1Len = 6   # "Their "
$full2 = $tr.Text
$tailStart = $p2Start + $run1Len

# Replace the tail run (previously non-underlined trailing clause) with the new
# closing phrase "imbalanced datasets".
$tailLen = $full2.Length - ($tailStart - 1)
$tr.Characters($tailStart, $tailLen).Text = "imbalanced datasets"

# Replace the leading run (previously "Their ") with the new opening sentence.
$tr.Characters($p2Start, $run1Len).Text = "Computation from contingency tables and performance characteristics under balanced and "

# The text box auto-fits its height to the text (<a:spAutoFit/>); the shorter new
# wording would otherwise shrink the shape, so restore the original height
# (584775 EMU = 46.04527559... pt) to keep the shape geometry unchanged.
$targetShape.Height = 584775 / 12700
